# Fix margins issue resulting from Jinja whitespace elimination.
# The template previously split certain Jinja tags across multiple Word
# runs (an artifact of manual formatting), which left the `-%}` tokens
# severed from their neighbours. Re-merge the run text so the intended
# Jinja whitespace-control tokens are correctly adjacent / removed.

$d = $word.ActiveDocument
$quoteOpen  = [char]8220
$quoteClose = [char]8221

# ---------------------------------------------------------------------
# Hunk 1: "{% for work in res.work %}" + "{% if not loop.first %}...{% endif -%}"
#         + (run with <w:br/> and work.name...) -> merge the first two
#         runs' text and let it coalesce with the following run.
# ---------------------------------------------------------------------
$searchText1 = "{% for work in res.work %}"
$found = $d.Content.Find.Execute($searchText1, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find hunk 1 anchor text"
}
$run1And2Text = $searchText1 + "{% if not loop.first %}{{ " + $quoteOpen + "\n" + $quoteClose + " }}{% endif -%}"

$combinedRange = $d.Content
$foundCombined = $combinedRange.Find.Execute($run1And2Text, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundCombined) {
    throw "Could not find hunk 1 combined run text"
}
$combinedRange.Delete()
$combinedRange.InsertAfter($run1And2Text)

# ---------------------------------------------------------------------
# Hunk 2: "{% " + "end" + "for %}" (three separate runs) -> single run
#         "{% endfor %}" (this is the paragraph immediately before the
#         "EDUCATION" paragraph).
# ---------------------------------------------------------------------
# Locate the target paragraph directly: it's the paragraph whose text
# renders as "{% endfor %}" and is immediately followed by the
# "EDUCATION" paragraph. (Paragraph.Range.Text includes the trailing
# paragraph-mark character, so trim it before comparing.)
$paras = $d.Paragraphs
$targetParaIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    $pText = $paras.Item($i).Range.Text.TrimEnd([char]13)
    if ($pText -eq "{% endfor %}") {
        if ($i -lt $paras.Count) {
            $nextText = $paras.Item($i + 1).Range.Text.TrimEnd([char]13)
            if ($nextText -eq "EDUCATION") {
                $targetParaIndex = $i
            }
        }
    }
}
if ($targetParaIndex -eq -1) {
    throw "Could not locate hunk 2 target paragraph"
}

$targetPara = $paras.Item($targetParaIndex)
$pStart = $targetPara.Range.Start
$pEnd = $targetPara.Range.End

# First run covers "{% " (3 characters); rewrite its text in place so it
# keeps its own run formatting (rPr), then delete the now-duplicated
# remainder of the paragraph (the old "end" + "for %}" runs).
$run1 = $d.Range($pStart, $pStart + 3)
if ($run1.Text -ne "{% ") {
    throw "Unexpected hunk 2 run1 text: $($run1.Text)"
}
$mergedText2 = "{% endfor %}"
$run1.Text = $mergedText2
$newRun1End = $pStart + $mergedText2.Length

$shift = $mergedText2.Length - 3
$restStart = $newRun1End
$restEnd = ($pEnd - 1) + $shift
if ($restEnd -gt $restStart) {
    $restRange = $d.Range($restStart, $restEnd)
    $restRange.Delete()
}

# ---------------------------------------------------------------------
# Hunk 3: final "{% endfor -%}" (last paragraph, right before sectPr)
#         -> "{% endfor %}"
# ---------------------------------------------------------------------
$paras2 = $d.Paragraphs
$lastPara = $paras2.Item($paras2.Count)
$lastParaText = $lastPara.Range.Text.TrimEnd([char]13)
if ($lastParaText -ne "{% endfor -%}") {
    throw "Unexpected last paragraph text: $lastParaText"
}
$lastRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$lastRange.Text = "{% endfor %}"

Write-Host "All three hunks applied."
